$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Create the new "2022-Q1" sheet, positioned right after "2021-Q4".
#    Duplicate "2021-Q4" wholesale (via Worksheet.Copy) so the new sheet
#    inherits identical sheet-level properties (sheetPr/outlinePr,
#    pageSetUpPr, pageMargins, column-A "index" style, header style, etc.)
#    instead of re-creating them by hand.
# ---------------------------------------------------------------------------
$q4 = $wb.Worksheets.Item("2021-Q4")
$q4.Copy($null, $q4)
$q1 = $wb.Worksheets.Item($q4.Index + 1)
$q1.Name = "2022-Q1"

# The template (2021-Q4) has 9 data rows (rows 2-10); 2022-Q1 only has 5
# data rows (rows 2-6), so drop the extra trailing rows.
$q1.Range("A7:H10").Delete()

# Columns B:G hold fund codes / decimal-looking figures that must stay text
# (e.g. "010460" would lose its leading zero, "3.41" would become a number)
# -- force them to literal text before writing, then strip the temporary
# "@" text format back off so the cells end up with no explicit style,
# matching the rest of the unstyled data cells.
$textCols = $q1.Range("B2:G6")
$textCols.NumberFormat = "@"

$q1.Range("B2").Value = "010460"
$q1.Range("C2").Value = "兴业研究精选混合"
$q1.Range("D2").Value = "3.41"
$q1.Range("E2").Value = "89.54"
$q1.Range("F2").Value = "3.01"
$q1.Range("G2").Value = "0.1026"
$q1.Range("H2").Value = 9

$q1.Range("B3").Value = "011603"
$q1.Range("C3").Value = "兴业高端制造混合A"
$q1.Range("D3").Value = "1.19"
$q1.Range("E3").Value = "76.40"
$q1.Range("F3").Value = "2.71"
$q1.Range("G3").Value = "0.0322"
$q1.Range("H3").Value = 8

$q1.Range("B4").Value = "011604"
$q1.Range("C4").Value = "兴业高端制造混合C"
$q1.Range("D4").Value = "0.54"
$q1.Range("E4").Value = "76.40"
$q1.Range("F4").Value = "2.71"
$q1.Range("G4").Value = "0.0146"
$q1.Range("H4").Value = 8

$q1.Range("B5").Value = "001866"
$q1.Range("C5").Value = "北信瑞丰新成长灵活配置混合"
$q1.Range("D5").Value = "0.07"
$q1.Range("E5").Value = "94.21"
$q1.Range("F5").Value = "9.04"
$q1.Range("G5").Value = "0.0063"
$q1.Range("H5").Value = 1

$q1.Range("B6").Value = "002303"
$q1.Range("C6").Value = "金鹰智慧生活灵活配置混合"
$q1.Range("D6").Value = "0.11"
$q1.Range("E6").Value = "89.88"
$q1.Range("F6").Value = "3.03"
$q1.Range("G6").Value = "0.0033"
$q1.Range("H6").Value = 10

$textCols.ClearFormats()

# ---------------------------------------------------------------------------
# 2. Update the "总计" (totals) sheet: add a new 2022-Q1 row, keeping the
#    existing 2021-Q4 row but pushed down to row 3.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Duplicate row 2 (2021-Q4, with its index-column style) down into row 3,
# then fix up its index value (0 -> 1).
$total.Range("A2:D2").Copy($total.Range("A3:D3"))
$total.Range("A3").Value = 1

# Overwrite row 2 in place with the new 2022-Q1 totals.
$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 5
$total.Range("D2").Value = 0.16

# ---------------------------------------------------------------------------
# 3. Restore the original active sheet ("2021-Q4") so the workbook-level
#    view state is unaffected by having created/renamed sheets above.
# ---------------------------------------------------------------------------
$q4.Activate()
